# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-level holding detail) right
#    before the "总计" (totals) sheet.
# 2) Insert a new leading row on the "总计" sheet summarising the new
#    quarter (date / holding count / holding market value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet, modelled on the existing per-quarter
#    sheets ("2021-Q4" supplies the header/column-A formatting to copy).
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalsBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalsBefore)
$newSheet.Name = "2022-Q1"
# Match the outline defaults used by the sibling per-quarter sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# NOTE: inserting a sheet shifts the position of everything from the
# insertion point onward, so a worksheet handle obtained *before* the
# Add() that referred to "总计" (or anything after it) is now stale -
# re-resolve it by name afterwards.
$totals = $wb.Worksheets.Item("总计")

# Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the header formatting (bold + border + centered) from the template sheet
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Fund-level detail rows
$fundRows = @(
    @{ B='213003'; C='宝盈策略增长混合'; D='10.28'; E='94.38'; F='8.41'; G='0.8645'; H=1 }
    @{ B='213002'; C='宝盈泛沿海增长混合'; D='5.10'; E='93.76'; F='9.64'; G='0.4916'; H=1 }
    @{ B='010330'; C='东吴兴享成长混合A'; D='11.63'; E='80.15'; F='3.78'; G='0.4396'; H=8 }
    @{ B='501201'; C='红土创新科技创新 3 年封闭运作灵活配置混合'; D='3.99'; E='96.70'; F='3.88'; G='0.1548'; H=6 }
    @{ B='011446'; C='长江新能源产业混合A'; D='2.61'; E='75.86'; F='4.68'; G='0.1221'; H=2 }
    @{ B='002707'; C='摩根士丹利华鑫科技领先灵活配置混合'; D='2.27'; E='93.05'; F='3.95'; G='0.0897'; H=9 }
    @{ B='970023'; C='天风天盈一年定期开放混合'; D='2.97'; E='39.47'; F='2.40'; G='0.0713'; H=9 }
    @{ B='001365'; C='大成正向回报灵活配置混合'; D='0.63'; E='75.86'; F='4.68'; G='0.0295'; H=2 }
    @{ B='011462'; C='东吴兴享成长混合C'; D='0.33'; E='80.15'; F='3.78'; G='0.0125'; H=8 }
    @{ B='002149'; C='嘉实新优选灵活配置混合'; D='0.22'; E='93.76'; F='3.97'; G='0.0087'; H=10 }
    @{ B='002862'; C='金信量化精选灵活配置混合'; D='0.16'; E='94.28'; F='4.08'; G='0.0065'; H=9 }
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = ($r - 2)
    # Fund code and the numeric-looking metrics are stored as *text* in the
    # source data (leading zeros, fixed decimal formatting) - a leading
    # apostrophe forces Excel to keep them as text instead of numbers.
    $newSheet.Cells.Item($r, 2).Value = "'" + $row.B
    $newSheet.Cells.Item($r, 3).Value = $row.C
    $newSheet.Cells.Item($r, 4).Value = "'" + $row.D
    $newSheet.Cells.Item($r, 5).Value = "'" + $row.E
    $newSheet.Cells.Item($r, 6).Value = "'" + $row.F
    $newSheet.Cells.Item($r, 7).Value = "'" + $row.G
    $newSheet.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Copy column-A formatting (bold + border + centered) from the template sheet
$template.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet data,
#    pushing the existing quarters down and renumbering column A.
# ---------------------------------------------------------------------
$totals.Rows.Item(2).Insert()
$totals.Range("A2:D2").ClearFormats()

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 11
$totals.Range("D2").Value = 2.29

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
